# RPA datasets push 2024-04-27
# Insert a new record for "하나33호스팩" above the existing "하나32호스팩"
# row (currently worksheet row 9), pushing that row and every row below it
# down by one. The other rows are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9, shifting rows 9-11 down to 10-12.
$ws.Rows.Item(9).Insert()

# The date-like columns (B, F, G) hold plain text such as "2024-04-15",
# not real Excel dates. Force Text format first so COM's auto-detection
# doesn't silently convert the literal strings into date serials.
$ws.Range("B9:G9").NumberFormat = "@"

# Populate the newly inserted row 9 with the "하나33호스팩" record.
$ws.Range("A9").Value = "하나"
$ws.Range("B9").Value = "2024-04-15"
$ws.Range("C9").Value = "하나33호스팩"
$ws.Range("D9").Value = "하나"
$ws.Range("E9").Value = "하나"
$ws.Range("F9").Value = "2024-04-18"
$ws.Range("G9").Value = "2024-04-24"
$ws.Range("H9").Value = 7000
$ws.Range("I9").Value = 3500000
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100

# Drop the temporary Text number format again so the new row ends up with
# the same default (unstyled) formatting as every other data row.
$ws.Range("B9:G9").ClearFormats()
